{"js": "// Replace each old answer/date string with its corresponding new value.\n// The mapping below was derived from the unified diff: every <w:t> run in\n// the document (the date heading plus all 100 arithmetic-answer table\n// cells) maps 1:1 from an old value to a new value, and every old value is\n// unique in the document, so a simple search-and-replace per pair is safe\n// and unambiguous.\nconst replacements = [\n  [\"2023-10-12 Thursday\", \"2023-10-13 Friday\"],\n  [\"19+29=48\", \"38-38=0\"],\n  [\"64-35=29\", \"85-63=22\"],\n  [\"93-19=74\", \"31-23=8\"],\n  [\"52-14=38\", \"41+6=47\"],\n  [\"79-30=49\", \"24+72=96\"],\n  [\"28+45=73\", \"79-6=73\"],\n  [\"64-18=46\", \"90-3=87\"],\n  [\"41-22=19\", \"78-46=32\"],\n  [\"70-49=21\", \"66-25=41\"],\n  [\"34-13=21\", \"44+21=65\"],\n  [\"81-31=50\", \"2+26=28\"],\n  [\"77-65=12\", \"65-54=11\"],\n  [\"95-15=80\", \"16+0=16\"],\n  [\"66-32=34\", \"94-91=3\"],\n  [\"88+0=88\", \"2+42=44\"],\n  [\"34+59=93\", \"99-30=69\"],\n  [\"84-67=17\", \"98-22=76\"],\n  [\"93-39=54\", \"47-34=13\"],\n  [\"93+3=96\", \"33+55=88\"],\n  [\"93-65=28\", \"39-20=19\"],\n  [\"8+15=23\", \"98-34=64\"],\n  [\"91-53=38\", \"19+42=61\"],\n  [\"7+65=72\", \"46-6=40\"],\n  [\"84-9=75\", \"66-39=27\"],\n  [\"75+21=96\", \"14+10=24\"],\n  [\"3+87=90\", \"59-21=38\"],\n  [\"44+20=64\", \"25+27=52\"],\n  [\"39+52=91\", \"8+61=69\"],\n  [\"86-80=6\", \"96-95=1\"],\n  [\"43+38=81\", \"63-54=9\"],\n  [\"67-48=19\", \"16+71=87\"],\n  [\"39-28=11\", \"48-45=3\"],\n  [\"10-6=4\", \"54-12=42\"],\n  [\"64-24=40\", \"67+29=96\"],\n  [\"41+16=57\", \"97-59=38\"],\n  [\"49+31=80\", \"72-63=9\"],\n  [\"98-25=73\", \"58+39=97\"],\n  [\"89-57=32\", \"86-68=18\"],\n  [\"37+48=85\", \"70-21=49\"],\n  [\"25+43=68\", \"65+17=82\"],\n  [\"1+41=42\", \"38-11=27\"],\n  [\"39+33=72\", \"53-36=17\"],\n  [\"25-22=3\", \"82-73=9\"],\n  [\"44+50=94\", \"29+48=77\"],\n  [\"82-30=52\", \"58-3=55\"],\n  [\"70+19=89\", \"16+23=39\"],\n  [\"18-2=16\", \"3+61=64\"],\n  [\"20-13=7\", \"67+19=86\"],\n  [\"27+32=59\", \"7+21=28\"],\n  [\"25+55=80\", \"84-48=36\"],\n  [\"5+62=67\", \"65+29=94\"],\n  [\"5+83=88\", \"35-30=5\"],\n  [\"24+15=39\", \"29+23=52\"],\n  [\"36-16=20\", \"20+5=25\"],\n  [\"55-7=48\", \"8+50=58\"],\n  [\"24-1=23\", \"78+6=84\"],\n  [\"47-43=4\", \"18+58=76\"],\n  [\"50+44=94\", \"38+53=91\"],\n  [\"10+45=55\", \"36+15=51\"],\n  [\"2+92=94\", \"62-55=7\"],\n  [\"79-20=59\", \"85-6=79\"],\n  [\"10+33=43\", \"97-81=16\"],\n  [\"93-60=33\", \"34+27=61\"],\n  [\"27+57=84\", \"4+10=14\"],\n  [\"0+36=36\", \"82-48=34\"],\n  [\"80-43=37\", \"31+8=39\"],\n  [\"63+32=95\", \"3+15=18\"],\n  [\"85-25=60\", \"45+46=91\"],\n  [\"79+11=90\", \"11+52=63\"],\n  [\"73+26=99\", \"44+46=90\"],\n  [\"2+4=6\", \"54-22=32\"],\n  [\"90-45=45\", \"64+20=84\"],\n  [\"55-37=18\", \"26+59=85\"],\n  [\"35+34=69\", \"50+2=52\"],\n  [\"7+71=78\", \"72-28=44\"],\n  [\"13+63=76\", \"14+64=78\"],\n  [\"22+10=32\", \"85-70=15\"],\n  [\"16+35=51\", \"36-12=24\"],\n  [\"33+44=77\", \"96-26=70\"],\n  [\"32+3=35\", \"11+59=70\"],\n  [\"76-32=44\", \"35+42=77\"],\n  [\"52+21=73\", \"47+7=54\"],\n  [\"35+37=72\", \"13+73=86\"],\n  [\"91-19=72\", \"34-31=3\"],\n  [\"52-43=9\", \"41-11=30\"],\n  [\"48-11=37\", \"73-10=63\"],\n  [\"34+51=85\", \"37+56=93\"],\n  [\"23+38=61\", \"92-48=44\"],\n  [\"76-31=45\", \"79-14=65\"],\n  [\"6+92=98\", \"2+73=75\"],\n  [\"40+52=92\", \"70+18=88\"],\n  [\"84-35=49\", \"94+0=94\"],\n  [\"97-56=41\", \"54+9=63\"],\n  [\"13+53=66\", \"84+2=86\"],\n  [\"69-25=44\", \"74+3=77\"],\n  [\"32+19=51\", \"45+7=52\"],\n  [\"12+86=98\", \"49+42=91\"],\n  [\"3+9=12\", \"30-25=5\"],\n  [\"35-12=23\", \"45-44=1\"],\n  [\"50-25=25\", \"94-57=37\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text to replace: \" + oldText);\n  }\n\n  // Each old value is unique in the document, so replace the first (only) hit.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each pair is (old answer/date text, new answer/date text). Every old\n# value is unique within the document, so Find/Replace-All per pair is\n# unambiguous -- it can only ever touch the single matching run.\n$pairs = @(\n    ,@(\"2023-10-12 Thursday\", \"2023-10-13 Friday\")\n    ,@(\"19+29=48\", \"38-38=0\")\n    ,@(\"64-35=29\", \"85-63=22\")\n    ,@(\"93-19=74\", \"31-23=8\")\n    ,@(\"52-14=38\", \"41+6=47\")\n    ,@(\"79-30=49\", \"24+72=96\")\n    ,@(\"28+45=73\", \"79-6=73\")\n    ,@(\"64-18=46\", \"90-3=87\")\n    ,@(\"41-22=19\", \"78-46=32\")\n    ,@(\"70-49=21\", \"66-25=41\")\n    ,@(\"34-13=21\", \"44+21=65\")\n    ,@(\"81-31=50\", \"2+26=28\")\n    ,@(\"77-65=12\", \"65-54=11\")\n    ,@(\"95-15=80\", \"16+0=16\")\n    ,@(\"66-32=34\", \"94-91=3\")\n    ,@(\"88+0=88\", \"2+42=44\")\n    ,@(\"34+59=93\", \"99-30=69\")\n    ,@(\"84-67=17\", \"98-22=76\")\n    ,@(\"93-39=54\", \"47-34=13\")\n    ,@(\"93+3=96\", \"33+55=88\")\n    ,@(\"93-65=28\", \"39-20=19\")\n    ,@(\"8+15=23\", \"98-34=64\")\n    ,@(\"91-53=38\", \"19+42=61\")\n    ,@(\"7+65=72\", \"46-6=40\")\n    ,@(\"84-9=75\", \"66-39=27\")\n    ,@(\"75+21=96\", \"14+10=24\")\n    ,@(\"3+87=90\", \"59-21=38\")\n    ,@(\"44+20=64\", \"25+27=52\")\n    ,@(\"39+52=91\", \"8+61=69\")\n    ,@(\"86-80=6\", \"96-95=1\")\n    ,@(\"43+38=81\", \"63-54=9\")\n    ,@(\"67-48=19\", \"16+71=87\")\n    ,@(\"39-28=11\", \"48-45=3\")\n    ,@(\"10-6=4\", \"54-12=42\")\n    ,@(\"64-24=40\", \"67+29=96\")\n    ,@(\"41+16=57\", \"97-59=38\")\n    ,@(\"49+31=80\", \"72-63=9\")\n    ,@(\"98-25=73\", \"58+39=97\")\n    ,@(\"89-57=32\", \"86-68=18\")\n    ,@(\"37+48=85\", \"70-21=49\")\n    ,@(\"25+43=68\", \"65+17=82\")\n    ,@(\"1+41=42\", \"38-11=27\")\n    ,@(\"39+33=72\", \"53-36=17\")\n    ,@(\"25-22=3\", \"82-73=9\")\n    ,@(\"44+50=94\", \"29+48=77\")\n    ,@(\"82-30=52\", \"58-3=55\")\n    ,@(\"70+19=89\", \"16+23=39\")\n    ,@(\"18-2=16\", \"3+61=64\")\n    ,@(\"20-13=7\", \"67+19=86\")\n    ,@(\"27+32=59\", \"7+21=28\")\n    ,@(\"25+55=80\", \"84-48=36\")\n    ,@(\"5+62=67\", \"65+29=94\")\n    ,@(\"5+83=88\", \"35-30=5\")\n    ,@(\"24+15=39\", \"29+23=52\")\n    ,@(\"36-16=20\", \"20+5=25\")\n    ,@(\"55-7=48\", \"8+50=58\")\n    ,@(\"24-1=23\", \"78+6=84\")\n    ,@(\"47-43=4\", \"18+58=76\")\n    ,@(\"50+44=94\", \"38+53=91\")\n    ,@(\"10+45=55\", \"36+15=51\")\n    ,@(\"2+92=94\", \"62-55=7\")\n    ,@(\"79-20=59\", \"85-6=79\")\n    ,@(\"10+33=43\", \"97-81=16\")\n    ,@(\"93-60=33\", \"34+27=61\")\n    ,@(\"27+57=84\", \"4+10=14\")\n    ,@(\"0+36=36\", \"82-48=34\")\n    ,@(\"80-43=37\", \"31+8=39\")\n    ,@(\"63+32=95\", \"3+15=18\")\n    ,@(\"85-25=60\", \"45+46=91\")\n    ,@(\"79+11=90\", \"11+52=63\")\n    ,@(\"73+26=99\", \"44+46=90\")\n    ,@(\"2+4=6\", \"54-22=32\")\n    ,@(\"90-45=45\", \"64+20=84\")\n    ,@(\"55-37=18\", \"26+59=85\")\n    ,@(\"35+34=69\", \"50+2=52\")\n    ,@(\"7+71=78\", \"72-28=44\")\n    ,@(\"13+63=76\", \"14+64=78\")\n    ,@(\"22+10=32\", \"85-70=15\")\n    ,@(\"16+35=51\", \"36-12=24\")\n    ,@(\"33+44=77\", \"96-26=70\")\n    ,@(\"32+3=35\", \"11+59=70\")\n    ,@(\"76-32=44\", \"35+42=77\")\n    ,@(\"52+21=73\", \"47+7=54\")\n    ,@(\"35+37=72\", \"13+73=86\")\n    ,@(\"91-19=72\", \"34-31=3\")\n    ,@(\"52-43=9\", \"41-11=30\")\n    ,@(\"48-11=37\", \"73-10=63\")\n    ,@(\"34+51=85\", \"37+56=93\")\n    ,@(\"23+38=61\", \"92-48=44\")\n    ,@(\"76-31=45\", \"79-14=65\")\n    ,@(\"6+92=98\", \"2+73=75\")\n    ,@(\"40+52=92\", \"70+18=88\")\n    ,@(\"84-35=49\", \"94+0=94\")\n    ,@(\"97-56=41\", \"54+9=63\")\n    ,@(\"13+53=66\", \"84+2=86\")\n    ,@(\"69-25=44\", \"74+3=77\")\n    ,@(\"32+19=51\", \"45+7=52\")\n    ,@(\"12+86=98\", \"49+42=91\")\n    ,@(\"3+9=12\", \"30-25=5\")\n    ,@(\"35-12=23\", \"45-44=1\")\n    ,@(\"50-25=25\", \"94-57=37\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format, ReplaceWith,\n    #         Replace(2=wdReplaceAll))\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find expected text to replace: $oldText\"\n    }\n}\n"}
